# Actualización automática 2025-06-27 17:25:45
$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet: VENTAS POR GRUPO
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("P18").Value = 187.29
$ws1.Range("Q18").Value = 44.6

$ws1.Range("E36").Value = 199.74
$ws1.Range("M36").Value = 73.34

$ws1.Range("M59").Value = 535.85

$ws1.Range("E149").Value = 129.64
$ws1.Range("G149").Value = 81.48

$ws1.Range("M154").Value = 2156.54

$ws1.Range("Q159").Value = 1935.16

$ws1.Range("H177").Value = 772.64
$ws1.Range("I177").Value = 739.5
$ws1.Range("M177").Value = 612.86
$ws1.Range("N177").Value = 2042.06

$ws1.Range("M257").Value = 25219.13

# Totals / counters row 279
$ws1.Range("E279").Value = "12 de 277"
$ws1.Range("G279").Value = "6 de 277"
$ws1.Range("H279").Value = "9 de 277"
$ws1.Range("I279").Value = "12 de 277"
$ws1.Range("M279").Value = "46 de 277"
$ws1.Range("N279").Value = "3 de 277"
$ws1.Range("Q279").Value = "10 de 277"

# ----------------------------------------------------------------------
# Sheet: VENTA MENSUAL
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("F18").Value = 1743.32
$ws2.Range("F36").Value = 1516.28
$ws2.Range("F59").Value = 1079.23
$ws2.Range("F149").Value = 211.12
$ws2.Range("F154").Value = 2156.54
$ws2.Range("F159").Value = 1935.16
$ws2.Range("F177").Value = 4167.06
$ws2.Range("F257").Value = 25219.13
$ws2.Range("F279").Value = 267272.46

# ----------------------------------------------------------------------
# Sheet: CUMPLIMIENTO MENSUAL
# ----------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Range("D10").Value = 228.89
$ws3.Range("E10").Value = 421.36
$ws3.Range("F10").Value = 0.3520030757400999

$ws3.Range("D14").Value = 44.6
$ws3.Range("E14").Value = 438.4
$ws3.Range("F14").Value = 0.09233954451345756

$ws3.Range("D19").Value = 5502.24
$ws3.Range("E19").Value = 317.7600000000002
$ws3.Range("F19").Value = 0.9454020618556701

$ws3.Range("D21").Value = 871.86
$ws3.Range("E21").Value = -225.86
$ws3.Range("F21").Value = 1.349628482972136

$ws3.Range("D32").Value = 11324.01
$ws3.Range("E32").Value = 4365.99
$ws3.Range("F32").Value = 0.721734225621415

$ws3.Range("D33").Value = 43477.13
$ws3.Range("E33").Value = 2268.559000000001
$ws3.Range("F33").Value = 0.9504093380252727

$ws3.Range("D70").Value = 129.64
$ws3.Range("E70").Value = 384.1910466593361
$ws3.Range("F70").Value = 0.2523008308720391

$ws3.Range("D72").Value = 125.02
$ws3.Range("E72").Value = -18.2
$ws3.Range("F72").Value = 1.170380078636959

$ws3.Range("D80").Value = 2366.8
$ws3.Range("E80").Value = -1883.8
$ws3.Range("F80").Value = 4.900207039337475

$ws3.Range("D82").Value = 5620.63
$ws3.Range("E82").Value = 23911.81
$ws3.Range("F82").Value = 0.1903205424272427

$ws3.Range("D90").Value = 772.64
$ws3.Range("E90").Value = 1327.36
$ws3.Range("F90").Value = 0.3679238095238095

$ws3.Range("D91").Value = 825.9
$ws3.Range("E91").Value = -75.89999999999998
$ws3.Range("F91").Value = 1.1012

$ws3.Range("D99").Value = 4371.46
$ws3.Range("E99").Value = 34045.71
$ws3.Range("F99").Value = 0.1137892249741457

$ws3.Range("D100").Value = 2268.66
$ws3.Range("E100").Value = -1926.66
$ws3.Range("F100").Value = 6.633508771929824

$ws3.Range("D119").Value = 4430.59
$ws3.Range("E119").Value = -4430.59

$ws3.Range("D120").Value = 594.78
$ws3.Range("E120").Value = 16905.22
$ws3.Range("F120").Value = 0.03398742857142857

$ws3.Range("D135").Value = 29659.77
$ws3.Range("E135").Value = -1449.93
$ws3.Range("F135").Value = 1.05139802281757

$ws3.Range("D138").Value = 292678.16
$ws3.Range("E138").Value = 129786.6105625342
$ws3.Range("F138").Value = 0.6927871396478422
